# Add test data for Delegate and ReplyAll transmittal actions.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before the existing "Action-Level3" column (O) to hold
# the new "DelegateTo" data; this shifts the old column O to column P and
# keeps all existing O4/O5 ("Submission") values moving along with it.
$ws.Columns("O").Insert()

# New header for the inserted column.
$ws.Range("O1").Value = "DelegateTo"

# Row 8: Delegate test data.
$ws.Range("A8").Value = "AutoTestAdmin"
$ws.Range("C8").Value = "New Transmittal from Automation"
$ws.Range("D8").Value = "UnTick"
$ws.Range("E8").Value = "Correspondence"
$ws.Range("F8").Value = "Request for Information"
$ws.Range("L8").Value = "Delegate- Message for New transmittal"
$ws.Range("M8").Value = "Delegate"
$ws.Range("O8").Value = "AutoTestUser"
$ws.Range("P8").Value = "Submission"

# Row 9: Reply All test data.
$ws.Range("A9").Value = "AutoTestAdmin"
$ws.Range("B9").Value = "AutoTestUser"
$ws.Range("C9").Value = "New Transmittal from Automation"
$ws.Range("D9").Value = "UnTick"
$ws.Range("E9").Value = "Correspondence"
$ws.Range("F9").Value = "Request for Information"
$ws.Range("L9").Value = "Reply All- Message for New transmittal"
$ws.Range("M9").Value = "ReplyAll"
$ws.Range("P9").Value = "Submission"

# Match the recorded selection left after entering the new data.
$ws.Range("F9").Select()

# Re-apply the sheet's "best fit" column widths now that new, wider content
# (e.g. the new DelegateTo/Message columns) has been added. Columns that did
# not receive any new/longer text (e.g. the relocated "Action-Level3" column)
# are left untouched so their existing best-fit width is preserved.
$ws.Columns("F").AutoFit()
$ws.Columns("L").AutoFit()
$ws.Columns("N").AutoFit()
$ws.Columns("O").AutoFit()
